# Update the "取得日時" (retrieved at) timestamp column for the captured
# rows on the "ランサーズ" sheet to reflect the latest scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-05 01:21:02"

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
